$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "67.336.28"
$ws.Cells.Item(2, 5).Value = "  +4.83%  "
$ws.Cells.Item(3, 4).Value = "3.238.38"
$ws.Cells.Item(3, 5).Value = "  +2.80%  "
$ws.Cells.Item(4, 5).Value = "  +0.00%  "
$ws.Cells.Item(5, 4).Value = "'576.70"
$ws.Cells.Item(5, 5).Value = "  +2.53%  "
$ws.Cells.Item(6, 4).Value = "'178.56"
$ws.Cells.Item(6, 5).Value = "  +6.32%  "
$ws.Cells.Item(7, 5).Value = "  -0.01%  "
$ws.Cells.Item(8, 4).Value = "'0.600"
$ws.Cells.Item(8, 5).Value = "  -2.44%  "
$ws.Cells.Item(9, 4).Value = "3.232.60"
$ws.Cells.Item(9, 5).Value = "  +2.70%  "
$ws.Cells.Item(10, 5).Value = "  +4.22%  "
$ws.Cells.Item(11, 5).Value = "  +3.32%  "
$ws.Cells.Item(12, 5).Value = "  +4.87%  "
$ws.Cells.Item(13, 4).Value = "3.794.59"
$ws.Cells.Item(13, 5).Value = "  +2.77%  "
$ws.Cells.Item(14, 5).Value = "  +0.99%  "
$ws.Cells.Item(15, 4).Value = "'27.77"
$ws.Cells.Item(15, 5).Value = "  +2.83%  "
$ws.Cells.Item(16, 4).Value = "67.233.82"
$ws.Cells.Item(16, 5).Value = "  +4.79%  "
$ws.Cells.Item(17, 5).Value = "  +2.85%  "
$ws.Cells.Item(18, 4).Value = "3.235.43"
$ws.Cells.Item(18, 5).Value = "  +2.96%  "
$ws.Cells.Item(19, 4).Value = "'5.79"
$ws.Cells.Item(19, 5).Value = "  +1.49%  "
$ws.Cells.Item(20, 5).Value = "  +3.65%  "
$ws.Cells.Item(21, 4).Value = "'374.06"
$ws.Cells.Item(21, 5).Value = "  +6.77%  "
$ws.Cells.Item(22, 4).Value = "'7.57"
$ws.Cells.Item(22, 5).Value = "  +5.60%  "
$ws.Cells.Item(23, 5).Value = "  +0.05%  "
$ws.Cells.Item(24, 4).Value = "'71.05"
$ws.Cells.Item(24, 5).Value = "  +4.53%  "
$ws.Cells.Item(25, 5).Value = "  +2.07%  "
$ws.Cells.Item(26, 2).Value = "PEPE"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(26, 4).Value = "'0.0000118"
$ws.Cells.Item(26, 5).Value = "  +3.20%  "
$ws.Cells.Item(27, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(27, 4).Value = "'9.60"
$ws.Cells.Item(27, 5).Value = "  +0.43%  "
$ws.Cells.Item(28, 2).Value = "Kaspa"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(28, 4).Value = "'0.181"
$ws.Cells.Item(28, 5).Value = "  +3.20%  "
$ws.Cells.Item(29, 2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Cells.Item(29, 4).Value = "'1.00"
$ws.Cells.Item(29, 5).Value = "  +0.42%  "
$ws.Cells.Item(30, 2).Value = "PancakeSwap"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(30, 4).Value = "'1.97"
$ws.Cells.Item(30, 5).Value = "  +4.51%  "
$ws.Cells.Item(31, 2).Value = "NEARProtocol"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(31, 4).Value = "'5.59"
$ws.Cells.Item(31, 5).Value = "  +2.92%  "
$ws.Cells.Item(32, 2).Value = "EthereumClassic"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(32, 4).Value = "'22.54"
$ws.Cells.Item(32, 5).Value = "  +3.22%  "
$ws.Cells.Item(33, 2).Value = "USDe"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(33, 4).Value = "'0.998"
$ws.Cells.Item(33, 5).Value = "  +0.00%  "
$ws.Cells.Item(34, 2).Value = "Fetch.AI"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(34, 4).Value = "'1.27"
$ws.Cells.Item(34, 5).Value = "  +6.86%  "
$ws.Cells.Item(35, 2).Value = "Aptos"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(35, 4).Value = "'6.79"
$ws.Cells.Item(35, 5).Value = "  +3.42%  "
$ws.Cells.Item(36, 2).Value = "Monero"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(36, 4).Value = "'163.95"
$ws.Cells.Item(36, 5).Value = "  +6.66%  "
$ws.Cells.Item(37, 2).Value = "ImmutableX"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(37, 4).Value = "'1.49"
$ws.Cells.Item(37, 5).Value = "  +4.88%  "
$ws.Cells.Item(38, 2).Value = "Mantle"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(38, 4).Value = "'0.859"
$ws.Cells.Item(38, 5).Value = "  +5.58%  "
$ws.Cells.Item(39, 2).Value = "Stacks"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(39, 4).Value = "'1.84"
$ws.Cells.Item(39, 5).Value = "  +8.99%  "
$ws.Cells.Item(40, 2).Value = "RenderToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(40, 4).Value = "'6.82"
$ws.Cells.Item(40, 5).Value = "  +14.80%  "
$ws.Cells.Item(41, 2).Value = "EnergySwap"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(41, 4).Value = "'26.61"
$ws.Cells.Item(41, 5).Value = "  +1.83%  "
$ws.Cells.Item(42, 2).Value = "Bittensor"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(42, 4).Value = "'362.60"
$ws.Cells.Item(42, 5).Value = "  +14.49%  "
$ws.Cells.Item(43, 2).Value = "dogwifhat"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(43, 4).Value = "'2.57"
$ws.Cells.Item(43, 5).Value = "  +5.41%  "
$ws.Cells.Item(44, 2).Value = "Maker"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(44, 4).Value = "2.708.26"
$ws.Cells.Item(44, 5).Value = "  +4.13%  "
$ws.Cells.Item(45, 2).Value = "Filecoin"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(45, 4).Value = "'4.39"
$ws.Cells.Item(45, 5).Value = "  +5.60%  "
$ws.Cells.Item(46, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(46, 4).Value = "'25.61"
$ws.Cells.Item(46, 5).Value = "  +8.04%  "
$ws.Cells.Item(47, 2).Value = "OKB"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(47, 4).Value = "'40.36"
$ws.Cells.Item(47, 5).Value = "  +2.85%  "
$ws.Cells.Item(48, 2).Value = "Hedera"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(48, 4).Value = "'0.0672"
$ws.Cells.Item(48, 5).Value = "  +4.32%  "
$ws.Cells.Item(49, 2).Value = "VeChain"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(49, 4).Value = "'0.0277"
$ws.Cells.Item(49, 5).Value = "  +2.54%  "
$ws.Cells.Item(50, 2).Value = "Stellar"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(50, 4).Value = "'0.103"
$ws.Cells.Item(50, 5).Value = "  +0.62%  "
$ws.Cells.Item(51, 2).Value = "ONDO"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Cells.Item(51, 4).Value = "'0.992"
$ws.Cells.Item(51, 5).Value = "  +6.14%  "

Write-Host "Applied 142 cell updates"
